# Update LDLC price history: a new snapshot column is inserted right
# before the "nom" / "url_produit" columns (currently FI / FJ), pushing
# them one column to the right (FJ / FK). The brand-new snapshot column
# gets header "2026-02-04 13:55:54" and, for every product row that
# already had a price in the previous last snapshot column (FH), the
# same price value is carried forward into the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FI is column index 165 (1-based). Inserting a whole column there
# shifts the existing FI ("nom") -> FJ and FJ ("url_produit") -> FK,
# exactly like the diff shows, and Excel grows the sheet dimension
# automatically (A1:FJ208 -> A1:FK208).
$ws.Range("FI1").EntireColumn.Insert()

# New header for the freshly inserted snapshot column.
$ws.Range("FI1").Value2 = "2026-02-04 13:55:54"

# Column FH (index 164) was the previous last price-snapshot column.
# Carry its value over into the new FI column (index 165) for every
# data row, mirroring what happened for every earlier snapshot column
# whenever a new price reading was appended. Rows where FH has no
# price (blank product rows) are left blank in FI as well.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $srcCell = $ws.Cells.Item($r, 164)   # FH<r>
    $dstCell = $ws.Cells.Item($r, 165)   # FI<r>
    $v = $srcCell.Value2
    if ($v -ne "") {
        $dstCell.Value2 = $v
    }
}
